# Update with restock suggestion
# Applies the diff to "Forecast Comparison" and "Summary" sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# --- Remove the old "Sales Volume Rank" column (Q). This shifts the former
# column R ("Lifecycle Stage") left into column Q, matching the new
# dimension A1:Q17 and the header move seen in the diff. ---
$ws.Range("Q1:Q17").Delete()

# --- Week_Start_Date (column B) now populated with the week's start date.
# A leading apostrophe forces these to be stored as text (matching the
# original inlineStr typing) rather than being auto-converted to date
# serial numbers. ---
$ws.Range("B2").Value = "'2025-02-02"
$ws.Range("B3").Value = "'2025-02-09"
$ws.Range("B4").Value = "'2025-02-16"
$ws.Range("B5").Value = "'2025-02-23"
$ws.Range("B6").Value = "'2025-03-02"
$ws.Range("B7").Value = "'2025-03-09"
$ws.Range("B8").Value = "'2025-03-16"
$ws.Range("B9").Value = "'2025-03-23"
$ws.Range("B10").Value = "'2025-03-30"
$ws.Range("B11").Value = "'2025-04-06"
$ws.Range("B12").Value = "'2025-04-13"
$ws.Range("B13").Value = "'2025-04-20"
$ws.Range("B14").Value = "'2025-04-27"
$ws.Range("B15").Value = "'2025-05-04"
$ws.Range("B16").Value = "'2025-05-11"
$ws.Range("B17").Value = "'2025-05-18"

# --- Sales Trend (column O): "Stable (-)" -> "Stable" ---
$ws.Range("O2").Value = "Stable"
$ws.Range("O3").Value = "Stable"
$ws.Range("O4").Value = "Stable"
$ws.Range("O5").Value = "Stable"
$ws.Range("O6").Value = "Stable"
$ws.Range("O7").Value = "Stable"
$ws.Range("O8").Value = "Stable"
$ws.Range("O9").Value = "Stable"
$ws.Range("O10").Value = "Stable"
$ws.Range("O11").Value = "Stable"
$ws.Range("O12").Value = "Stable"
$ws.Range("O13").Value = "Stable"
$ws.Range("O14").Value = "Stable"
$ws.Range("O15").Value = "Stable"
$ws.Range("O16").Value = "Stable"
$ws.Range("O17").Value = "Stable"

# --- Seasonality Index (column P) updated values ---
$ws.Range("P2").Value = 1.07
$ws.Range("P3").Value = 1.01
$ws.Range("P4").Value = 0.81
$ws.Range("P5").Value = 1.08
$ws.Range("P6").Value = 0.85
$ws.Range("P7").Value = 1.08
$ws.Range("P8").Value = 0.85
$ws.Range("P9").Value = 1.19
$ws.Range("P10").Value = 0.92
$ws.Range("P11").Value = 0.81
$ws.Range("P12").Value = 0.9
$ws.Range("P13").Value = 0.85
$ws.Range("P14").Value = 0.85
$ws.Range("P15").Value = 0.97
$ws.Range("P16").Value = 1.05
$ws.Range("P17").Value = 0.92

# --- Lifecycle Stage (now column Q, after old Q "Sales Volume Rank" removed
# and R shifted in). All weeks move from "Growth" to "Mature". ---
$ws.Range("Q2").Value = "Mature"
$ws.Range("Q3").Value = "Mature"
$ws.Range("Q4").Value = "Mature"
$ws.Range("Q5").Value = "Mature"
$ws.Range("Q6").Value = "Mature"
$ws.Range("Q7").Value = "Mature"
$ws.Range("Q8").Value = "Mature"
$ws.Range("Q9").Value = "Mature"
$ws.Range("Q10").Value = "Mature"
$ws.Range("Q11").Value = "Mature"
$ws.Range("Q12").Value = "Mature"
$ws.Range("Q13").Value = "Mature"
$ws.Range("Q14").Value = "Mature"
$ws.Range("Q15").Value = "Mature"
$ws.Range("Q16").Value = "Mature"
$ws.Range("Q17").Value = "Mature"

# --- Inventory Coverage (column L) newly populated for a few rows ---
$ws.Range("L13").Value = 7.14
$ws.Range("L14").Value = 6.14
$ws.Range("L17").Value = 5.14

# --- Summary sheet: Max/Min Forecast Week become "N/A" ---
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Range("B13").Value = "N/A"
$ws2.Range("B15").Value = "N/A"
